$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Change Management Overview" ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

# Project name change
$ws1.Range("B6").Value = "Enterprise Cloud Infrastructure Migration"

# Insert an empty row marker at row 13 (between row 12 and row 14)
$ws1.Rows.Item(13).OutlineLevel = 0

# Objective text updates (AI/ML -> IT)
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new IT systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in IT technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for IT transformation"

# Insert an empty row marker at row 21 (between row 20 and row 22)
$ws1.Rows.Item(21).OutlineLevel = 0

# --- Sheet 2: "Change Impact Assessment" ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

# Insert an empty row marker at row 2 (between row 1 and row 3)
$ws2.Rows.Item(2).OutlineLevel = 0

# Stakeholder group / content updates
$ws2.Range("A4").Value = "IT Managers"
$ws2.Range("G4").Value = "IT automation"
$ws2.Range("A5").Value = "System Administrators"

# --- Sheet 3: "Change Activities" ---
$ws3 = $wb.Worksheets.Item("Change Activities")

# Insert an empty row marker at row 2 (between row 1 and row 3)
$ws3.Rows.Item(2).OutlineLevel = 0

Write-Host "Edits applied"
